$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (engine quantizes ColumnWidth to 1/6-character steps, so these
#     are chosen to land on the closest achievable stored width to the target) ---
$ws.Columns.Item(1).ColumnWidth = 22.333333333333332
$ws.Columns.Item(2).ColumnWidth = 19.666666666666668
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668
$ws.Columns.Item(4).ColumnWidth = 22.5
$ws.Columns.Item(5).ColumnWidth = 18.833333333333332
$ws.Columns.Item(6).ColumnWidth = 23.0
$ws.Columns.Item(7).ColumnWidth = 25.0

# --- Updated data values (registration results, incl. new mask-registration run) ---
$ws.Range("B2").Value = 0.38240000000000002
$ws.Range("C2").Value = 0.20757999999999999
$ws.Range("D2").Value = 0.81938999999999995
$ws.Range("E2").Value = 0.64339999999999997
$ws.Range("F2").Value = 0.84104000000000001
$ws.Range("G2").Value = 0.69471000000000005

$ws.Range("B3").Value = 0.049489999999999999
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.85614000000000001
$ws.Range("E3").Value = 0.84455999999999998
$ws.Range("F3").Value = 0.85511999999999999
$ws.Range("G3").Value = 0.84243000000000001

$ws.Range("B4").Value = 0.30375999999999997
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.88461000000000001
$ws.Range("E4").Value = 0.76642999999999994
$ws.Range("F4").Value = 0.88099000000000005
$ws.Range("G4").Value = 0.72480999999999995

$ws.Range("B5").Value = 0.40884999999999999
$ws.Range("C5").Value = 0.27077000000000001
$ws.Range("D5").Value = 0.93189
$ws.Range("E5").Value = 0.94987999999999995
$ws.Range("F5").Value = 0.93078000000000005
$ws.Range("G5").Value = 0.94355

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.87812000000000001
$ws.Range("E6").Value = 0.87246000000000001
$ws.Range("F6").Value = 0.872
$ws.Range("G6").Value = 0.86456

$ws.Range("B7").Value = 0.036530000000000002
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.56352000000000002
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0.56145999999999996
$ws.Range("G7").Value = 0

$ws.Range("B8").Value = 0.38774999999999998
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0.74173999999999995
$ws.Range("E8").Value = 0.70555000000000001
$ws.Range("F8").Value = 0.73912999999999995
$ws.Range("G8").Value = 0.68742000000000003

$ws.Range("B9").Value = 0.24389
$ws.Range("C9").Value = 0.14654
$ws.Range("D9").Value = 0.66715999999999998
$ws.Range("E9").Value = 0.64102000000000003
$ws.Range("F9").Value = 0.66966000000000003
$ws.Range("G9").Value = 0.64731000000000005

$ws.Range("B10").Value = 0.46000999999999997
$ws.Range("C10").Value = 0.31574000000000002
$ws.Range("D10").Value = 0.82872999999999997
$ws.Range("E10").Value = 0.77553000000000005
$ws.Range("F10").Value = 0.83130999999999999
$ws.Range("G10").Value = 0.77678999999999998

$ws.Range("B11").Value = 0.14257
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.78188999999999997
$ws.Range("E11").Value = 0.61377000000000004
$ws.Range("F11").Value = 0.79803000000000002
$ws.Range("G11").Value = 0.62234

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0.56352000000000002
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0.56145999999999996
$ws.Range("G12").Value = 0

$ws.Range("B13").Value = 0.46000999999999997
$ws.Range("C13").Value = 0.31574000000000002
$ws.Range("D13").Value = 0.93189
$ws.Range("E13").Value = 0.94987999999999995
$ws.Range("F13").Value = 0.93078000000000005
$ws.Range("G13").Value = 0.94355

$ws.Range("B14").Value = 0.23960500000000001
$ws.Range("C14").Value = 0.1046975
$ws.Range("D14").Value = 0.78738333333333321
$ws.Range("E14").Value = 0.6468733333333333
$ws.Range("F14").Value = 0.78931333333333331
$ws.Range("G14").Value = 0.6456225000000001

$ws.Range("B15").Value = 0.24389
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0.81938999999999995
$ws.Range("E15").Value = 0.70555000000000001
$ws.Range("F15").Value = 0.83130999999999999
$ws.Range("G15").Value = 0.69471000000000005

# --- Highlight the F14 cell (new mask-registration result) with a bold font + yellow fill ---
$f14 = $ws.Range("F14")
$f14.Font.Bold = $true
$f14.Interior.Color = 65535

# --- Selection moves to F14 ---
$f14.Select() | Out-Null
